$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of I2 and J2 (keep formatting/style)
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Update the selection to H2:O11
$ws.Range("H2:O11").Select()
